$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 643
$ws.Range("I19").Value = 501.8
$ws.Range("J19").Value = 784.2
$ws.Range("K19").Value = 501.8
$ws.Range("L19").Value = 784.2
$ws.Range("M19").Value = -326.8
$ws.Range("N19").Value = -1134.2
$ws.Range("H123").Value = 98796
$ws.Range("J123").Value = 98796
$ws.Range("L123").Value = 98796
$ws.Range("N123").Value = -108596
$ws.Range("H129").Value = 938.7143
$ws.Range("I129").Value = 533.1111
$ws.Range("J129").Value = 1006.3148
$ws.Range("K129").Value = 1599.3333
$ws.Range("L129").Value = 3018.9444
$ws.Range("M129").Value = 3400.6667
$ws.Range("N129").Value = -13018.9444
$ws.Range("H132").Value = 292938.9
$ws.Range("I132").Value = 338371.3
$ws.Range("J132").Value = 59286.57
$ws.Range("K132").Value = 1015113.9
$ws.Range("L132").Value = 177859.71
$ws.Range("M132").Value = -1012583.9
$ws.Range("N132").Value = -182919.71
$ws.Range("H137").Value = 34484348
$ws.Range("I137").Value = 40001452
$ws.Range("J137").Value = 2451
$ws.Range("K137").Value = 120004356
$ws.Range("L137").Value = 7353
$ws.Range("M137").Value = -120001806
$ws.Range("N137").Value = -12453
$ws.Range("H138").Value = 5433369.5
$ws.Range("I138").Value = 1076033
$ws.Range("K138").Value = 3228099
$ws.Range("M138").Value = -3222959

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18242.984
$ws.Range("I32").Value = 2177.3777
$ws.Range("J32").Value = 52669.285
$ws.Range("K32").Value = 2177.3777
$ws.Range("L32").Value = 52669.285
$ws.Range("M32").Value = -1890.3777
$ws.Range("N32").Value = -53243.285
$ws.Range("H45").Value = 1123.5555
$ws.Range("I45").Value = 1030.2858
$ws.Range("K45").Value = 1030.2858
$ws.Range("M45").Value = -653.2858000000001
$ws.Range("H110").Value = 858.36365
$ws.Range("I110").Value = 815.7778
$ws.Range("J110").Value = 1050
$ws.Range("K110").Value = 815.7778
$ws.Range("L110").Value = 1050
$ws.Range("M110").Value = 1229.2222
$ws.Range("N110").Value = -5140
$ws.Range("H122").Value = 1195.8667
$ws.Range("I122").Value = 999
$ws.Range("J122").Value = 1420.8572
$ws.Range("K122").Value = 2997
$ws.Range("L122").Value = 4262.571599999999
$ws.Range("M122").Value = -547
$ws.Range("N122").Value = -9162.571599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1351.0714
$ws.Range("I20").Value = 1138.5264
$ws.Range("J20").Value = 1799.7778
$ws.Range("K20").Value = 1138.5264
$ws.Range("L20").Value = 1799.7778
$ws.Range("M20").Value = -891.5264
$ws.Range("N20").Value = -2293.7778
$ws.Range("H107").Value = 621.2917
$ws.Range("I107").Value = 571.82355
$ws.Range("J107").Value = 741.4286
$ws.Range("K107").Value = 571.82355
$ws.Range("L107").Value = 741.4286
$ws.Range("M107").Value = 1348.17645
$ws.Range("N107").Value = -4581.4286
$ws.Range("H134").Value = 3414.0908
$ws.Range("I134").Value = 2069.375
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 6208.125
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -3673.125
$ws.Range("N134").Value = -26070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1953.2273
$ws.Range("I58").Value = 847.3570999999999
$ws.Range("J58").Value = 3888.5
$ws.Range("K58").Value = 847.3570999999999
$ws.Range("L58").Value = 3888.5
$ws.Range("M58").Value = -644.3570999999999
$ws.Range("N58").Value = -4294.5
$ws.Range("H122").Value = 2408.077
$ws.Range("I122").Value = 1391.3636
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 4174.0908
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -1724.0908
$ws.Range("N122").Value = -28900
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
$ws.Range("H136").Value = 1953.2273
$ws.Range("I136").Value = 847.3570999999999
$ws.Range("J136").Value = 3888.5
$ws.Range("K136").Value = 2542.0713
$ws.Range("L136").Value = 11665.5
$ws.Range("M136").Value = 7.92870000000039
$ws.Range("N136").Value = -16765.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 629.8125
$ws.Range("I121").Value = 95.8
$ws.Range("J121").Value = 872.5454999999999
$ws.Range("K121").Value = 287.4
$ws.Range("L121").Value = 2617.6365
$ws.Range("M121").Value = 1022.6
$ws.Range("N121").Value = -5237.6365
$ws.Range("H122").Value = 911.3333
$ws.Range("I122").Value = 281.2
$ws.Range("J122").Value = 1153.6923
$ws.Range("K122").Value = 2530.8
$ws.Range("L122").Value = 10383.2307
$ws.Range("M122").Value = -80.79999999999973
$ws.Range("N122").Value = -15283.2307
$ws.Range("H129").Value = 2202.6
$ws.Range("I129").Value = 644.5
$ws.Range("J129").Value = 2592.125
$ws.Range("K129").Value = 1933.5
$ws.Range("L129").Value = 7776.375
$ws.Range("M129").Value = 3066.5
$ws.Range("N129").Value = -17776.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 652
$ws.Range("I5").Value = 304
$ws.Range("K5").Value = 304
$ws.Range("M5").Value = -192

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H22").Value = 16925.334
$ws.Range("I22").Value = 337.5
$ws.Range("J22").Value = 50101
$ws.Range("K22").Value = 337.5
$ws.Range("L22").Value = 50101
$ws.Range("M22").Value = -42.5
$ws.Range("N22").Value = -50691
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H27").Value = 16925.334
$ws.Range("I27").Value = 337.5
$ws.Range("J27").Value = 50101
$ws.Range("K27").Value = 337.5
$ws.Range("L27").Value = 50101
$ws.Range("M27").Value = -230.5
$ws.Range("N27").Value = -50315
$ws.Range("H93").Value = 633.35
$ws.Range("I93").Value = 630.25
$ws.Range("J93").Value = 645.75
$ws.Range("K93").Value = 630.25
$ws.Range("L93").Value = 645.75
$ws.Range("M93").Value = 617.75
$ws.Range("N93").Value = -3141.75
$ws.Range("H122").Value = 4000.3333
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4000.3572
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 12001.0716
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -16901.0716

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 5500
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 1000
$ws.Range("M20").Value = -760
$ws.Range("H113").Value = 450
$ws.Range("J113").Value = 600
$ws.Range("L113").Value = 1800
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 127187.875
$ws.Range("I122").Value = 252126
$ws.Range("J122").Value = 2249.75
$ws.Range("K122").Value = 756378
$ws.Range("L122").Value = 6749.25
$ws.Range("M122").Value = -753928
$ws.Range("N122").Value = -11649.25
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 78943
$ws.Range("J141").Value = 78943
$ws.Range("L141").Value = 78943
$ws.Range("N141").Value = -89303
